$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.375.72'
$ws.Range("E2").Value = '  -0.26%  '
$ws.Range("D3").Value = '2.649.49'
$ws.Range("E3").Value = '  +2.57%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("E5").Value = '  +2.77%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.96'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.19%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.585'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.44%  '
$ws.Range("D9").Value = '2.648.91'
$ws.Range("E9").Value = '  +2.64%  '
$ws.Range("E10").Value = '  +1.88%  '
$ws.Range("E11").Value = '  +0.36%  '
$ws.Range("E12").Value = '  +0.36%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.364'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.44'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.34%  '
$ws.Range("D15").Value = '3.124.15'
$ws.Range("E15").Value = '  +2.48%  '
$ws.Range("D16").Value = '63.181.16'
$ws.Range("E16").Value = '  -0.39%  '
$ws.Range("E17").Value = '  -0.43%  '
$ws.Range("D18").Value = '2.663.94'
$ws.Range("E18").Value = '  +3.24%  '
$ws.Range("E19").Value = '  +2.52%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.46'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.24%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '342.80'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("E22").Value = '  +3.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '67.08'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.94%  '
$ws.Range("E25").Value = '  +2.47%  '
$ws.Range("E26").Value = '  -2.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.70'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.35%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '543.65'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +15.47%  '
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.95'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.14%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.05'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.83%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.81'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.79%  '
$ws.Range("D34").Value = '0.0₃0812'
$ws.Range("E34").Value = '  +1.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '171.66'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.82%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.16'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +14.03%  '
$ws.Range("E37").Value = '  +1.43%  '
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.14'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.14%  '
$ws.Range("E40").Value = '  +6.68%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '172.30'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.31%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("E43").Value = '  +1.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.46'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.75%  '
$ws.Range("E45").Value = '  +7.80%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.633'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.23%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0963'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.14%  '
$ws.Range("E48").Value = '  +1.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.86'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.30%  '
$ws.Range("E50").Value = '  +2.96%  '
$ws.Range("E51").Value = '  -0.90%  '
